# Apply the "core functional test cases of company" update:
#  - Re-word several Test Objective (column C) cells to reflect the
#    re-ordered / re-scoped company test suite.
#  - Row 2 (TC_001) and rows 3-4 (TC_002/TC_003) no longer have a final
#    PASSED/FAILED outcome recorded - they revert to "Test not executed" /
#    "Not Run" (the plain, unfilled status style).
#  - Row 14 (TC_013) now records a FAILED outcome (red fill).
#  - The three trailing rows (TC_014/TC_015/TC_016) are removed, shrinking
#    the sheet from 17 to 14 used rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - reuse an existing "status style" (fill/border/alignment)
# instead of Excel inventing a brand-new cell format for every touched cell.
$xlPasteFormats = -4122

# --- 1. Re-word the Test Objective column (C) -----------------------------

$ws.Range("C3").Value  = "Verify all mandatory field validation errors are displayed when form is submitted empty."
$ws.Range("C4").Value  = "Verify company name minimum length validation error message."
$ws.Range("C5").Value  = "Verify company name maximum length validation error message."
$ws.Range("C6").Value  = "Verify company name special character validation error message."
$ws.Range("C7").Value  = "Verify file upload size validation error."
$ws.Range("C8").Value  = "Verify file type upload validation error."
$ws.Range("C9").Value  = "Verify company creation with all optional fields including image upload."
$ws.Range("C10").Value = "Verify successful company creation with all mandatory fields."
$ws.Range("C11").Value = "Verify duplicate company name validation using the company created in TC_09."
$ws.Range("C12").Value = "Verify navigation to company details page via clicking first company."
$ws.Range("C13").Value = "Verify deletion of company created in TC_09."
$ws.Range("C14").Value = "Verify bulk deletion of few individual companies (3-5) using individual checkboxes."

# --- 2. Row 14 (TC_013): Not Run -> FAILED ---------------------------------
#    (done before rows 3/4 are converted away from their FAILED styling, so
#    there is still a FAILED-styled cell available to copy the format from)

$ws.Range("H3").Copy()
$ws.Range("H14").PasteSpecial($xlPasteFormats)
$ws.Range("I3").Copy()
$ws.Range("I14").PasteSpecial($xlPasteFormats)
$ws.Range("H14").Value = "Test failed - actual behavior did not match expected result"
$ws.Range("I14").Value = "FAILED"

# --- 3. Row 2 (TC_001): PASSED -> Not Run ----------------------------------

$ws.Range("H5").Copy()
$ws.Range("H2").PasteSpecial($xlPasteFormats)
$ws.Range("I5").Copy()
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("H2").Value = "Test not executed"
$ws.Range("I2").Value = "Not Run"

# --- 4. Rows 3-4 (TC_002/TC_003): FAILED -> Not Run ------------------------

$ws.Range("H5").Copy()
$ws.Range("H3").PasteSpecial($xlPasteFormats)
$ws.Range("I5").Copy()
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("H3").Value = "Test not executed"
$ws.Range("I3").Value = "Not Run"

$ws.Range("H5").Copy()
$ws.Range("H4").PasteSpecial($xlPasteFormats)
$ws.Range("I5").Copy()
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("H4").Value = "Test not executed"
$ws.Range("I4").Value = "Not Run"

# --- 5. Drop the trailing TC_014 / TC_015 / TC_016 rows --------------------

$ws.Rows("15:17").Delete()

$excel.CutCopyMode = $false
